$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Apply base formatting to the two new rows (8 and 9) first, mirroring
#     the style used by the existing data rows: wrap text + top vertical
#     alignment for B:Q, and top-alignment only (no wrap) for the SL column A.
$ws.Range("A8:Q9").VerticalAlignment = -4160
$ws.Range("A8:Q9").WrapText = $true
$ws.Range("A8").WrapText = $false
$ws.Range("A9").WrapText = $false

# --- Row 7 (SL 7): "Combo list selection is not working (ref: SEC_1.xml)"
$ws.Range("F8").Value = "Chosen option should be saved properly"
$ws.Range("G8").Value = "No option is showed as chosen"
$ws.Range("E8").Value = "1. Load any form in mobile                     2. Select any option from the option list                                                                   3. Click 'Save and Exit'    "

# --- Row 8 (SL 8): "Modification not works during synchronization (ref: facility profile)"
$ws.Range("F9").Value = "Updated data should be showed properly during sychronization in MS Access"
$ws.Range("G9").Value = "Updated data are not showing properly. "
$ws.Range("E9").Value = "1. Load any facility profile in mobile                                                 2. Do any modification in text field and option list                                            3. Click 'Send Now'                                   4. Sync MS Access with Cloud         "

# --- Titles added last, matching authoring order captured in the shared
#     string table.
$ws.Range("B8").Value = "Combo list selection is not working (ref: SEC_1.xml)"
$ws.Range("B9").Value = "Modification not works during synchronization (ref: facility profile)"

# --- Remaining columns for both rows.
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("I8").Value = "High"
$ws.Range("I9").Value = "High"

# --- Row heights to match the authored layout.
$ws.Range("A8").EntireRow.RowHeight = 60
$ws.Range("A9").EntireRow.RowHeight = 75

# --- View state: scroll/selection moved down to the newly-added rows.
[void]$ws.Range("D6").Select()
